$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text "Ready for handoff" -> "In Translation" wherever it
# appears (Overview!E2:F2, zh-cn!C2, de-de!C2). Cells.Replace mutates the
# shared-string text in place instead of dropping/re-adding a new string.
$null = $overview.Cells.Replace("Ready for handoff", "In Translation")
$null = $zhcn.Cells.Replace("Ready for handoff", "In Translation")
$null = $dede.Cells.Replace("Ready for handoff", "In Translation")

# Narrow the Status-related columns (was ~17.22 characters, now ~13.41).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
